$d = $word.ActiveDocument

# 1) Shorten the "processos logísticos ..." sentence in the FINALIDADE paragraph.
$rng = $d.Content
$rng.Find.Execute(
    "rocessos logísticos para tratar de atividades envolvendo setores da NNAQ",
    $true, $false, $false, $false, $false, $true, 0, $false,
    "rocessos logísticos", 2) | Out-Null

# 2) Expand the ÂMBITO sentence to add NCCP / DACI / Direção do CELOG.
$rng = $d.Content
$rng.Find.Execute(
    ", de observância obrigatória, aplica-se à Divisão de Nacionalização e Qualificação (NNAQ) do Centro Logístico da Aeronáutica (CELOG).",
    $true, $false, $false, $false, $false, $true, 0, $false,
    ", de observância obrigatória, aplica-se à Divisão de Nacionalização e Qualificação (NNAQ), à Seção de Controle de Publicação (NCCP), à Assessoria de Controle Interno (DACI) e à Direção do Centro Logístico da Aeronáutica (CELOG).", 2) | Out-Null

# 3) Rename heading "PROCESSOS SUPERIORES" -> "MACROPROCESSO".
$rng = $d.Content
$rng.Find.Execute(
    "PROCESSOS SUPERIORES",
    $true, $false, $false, $false, $false, $true, 0, $false,
    "MACROPROCESSO", 2) | Out-Null

# 4) Rename heading "PROCESSOS SUBORDINADOS" -> "SUBPROCESSOS".
$rng = $d.Content
$rng.Find.Execute(
    "PROCESSOS SUBORDINADOS",
    $true, $false, $false, $false, $false, $true, 0, $false,
    "SUBPROCESSOS", 2) | Out-Null

# 5) Complete the signature line "Chefe da Divisão de Nacionalização" -> "... e Qualificação".
$rng = $d.Content
$rng.Find.Execute(
    "Chefe da Divisão de Nacionalização",
    $true, $false, $false, $false, $false, $true, 0, $false,
    "Chefe da Divisão de Nacionalização e Qualificação", 2) | Out-Null

# 6) Update the cached PAGE field value in the footer from "2" to "1".
#    (Find/Replace does not reliably commit edits inside HeaderFooter ranges
#    in this runtime, so walk the footer's raw Characters collection and
#    rewrite the single "2" character directly - this keeps the PAGE/
#    NUMPAGES field codes and run formatting untouched.)
$ftr = $d.Sections.Item(1).Footers.Item(1)
if ($ftr.Exists) {
    $chars = $ftr.Range.Characters
    for ($i = 1; $i -le $chars.Count; $i++) {
        $ch = $chars.Item($i)
        if ($ch.Text -eq "2") {
            $ch.Text = "1"
            break
        }
    }
}
